$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert()
$ws.Activate()
$ws.Range("R4").Select()
